$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------
# Edit 1: paragraph 1 - append a new run "的撒法发是否" after the
# existing run, inside the same paragraph.
# ---------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1TextRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)  # exclude paragraph mark

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:wordDocument ' + $wNs + '><w:body><w:p>'
$xml1 += '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>是的规范化买奶粉大热国会女Greg NBC个人</w:t></w:r>'
$xml1 += '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的撒法发是否</w:t></w:r>'
$xml1 += '</w:p></w:body></w:wordDocument>'
$p1TextRange.InsertXML($xml1)

# ---------------------------------------------------------------
# Edit 2: last paragraph - split the single run into three runs
# separated by proofErr gramStart/gramEnd markers, and drop the
# paragraph-mark run formatting (pPr/rPr/rFonts hint=eastAsia).
# ---------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$start = $pLast.Range.Start
$end = $pLast.Range.End   # includes the paragraph mark

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:wordDocument ' + $wNs + '><w:body>'
$xml2 += '<w:p></w:p>'
$xml2 += '<w:p>'
$xml2 += '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>大厦的功法哈佛嘎</w:t></w:r>'
$xml2 += '<w:proofErr w:type="gramStart"/>'
$xml2 += '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>哈好的</w:t></w:r>'
$xml2 += '<w:proofErr w:type="gramEnd"/>'
$xml2 += '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>对</w:t></w:r>'
$xml2 += '</w:p>'
$xml2 += '</w:body></w:wordDocument>'

$lastRange = $d.Range($start, $end)
$lastRange.InsertXML($xml2)

# The insertion above leaves a spurious empty paragraph mark right
# before our freshly-inserted paragraph (since the original
# paragraph mark gets preserved as a separate, now-empty
# paragraph). Delete that leftover mark to merge it away.
$extraMark = $d.Range($start, $start + 1)
$extraMark.Delete()
